# Reorders the worker arrears (Estado de Cuenta) rows so each worker's
# periods are grouped together (most recent period 2108 first, oldest
# 2102 last), and updates KEVIN ANGULO MARIMON's "Salario Basico"
# (base salary) from 5,000,000 to 1,250,000 for every one of his rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data block (rows 16-36, cols B-J) before rewriting it
$ws.Range("B16:J36").ClearContents()

# JUAN CARLOS MATOS CAVADIA
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73582320"
$ws.Range("D16").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E16").Value = "2108"
$ws.Range("F16").Value = 35129
$ws.Range("G16").Value = 1817052

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73582320"
$ws.Range("D17").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E17").Value = "2107"
$ws.Range("F17").Value = 72682
$ws.Range("G17").Value = 1817052

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73582320"
$ws.Range("D18").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E18").Value = "2106"
$ws.Range("F18").Value = 72682
$ws.Range("G18").Value = 1817052

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73582320"
$ws.Range("D19").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E19").Value = "2105"
$ws.Range("F19").Value = 72682
$ws.Range("G19").Value = 1817052

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73582320"
$ws.Range("D20").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E20").Value = "2104"
$ws.Range("F20").Value = 72682
$ws.Range("G20").Value = 1817052

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73582320"
$ws.Range("D21").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E21").Value = "2103"
$ws.Range("F21").Value = 72682
$ws.Range("G21").Value = 1817052

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73582320"
$ws.Range("D22").Value = "JUAN CARLOS MATOS CAVADIA"
$ws.Range("E22").Value = "2102"
$ws.Range("F22").Value = 72682
$ws.Range("G22").Value = 1817052

# KEVIN ANGULO MARIMON
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143357213"
$ws.Range("D23").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E23").Value = "2108"
$ws.Range("F23").Value = 35129
$ws.Range("G23").Value = 1250000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143357213"
$ws.Range("D24").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E24").Value = "2107"
$ws.Range("F24").Value = 200000
$ws.Range("G24").Value = 1250000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1143357213"
$ws.Range("D25").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E25").Value = "2106"
$ws.Range("F25").Value = 200000
$ws.Range("G25").Value = 1250000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1143357213"
$ws.Range("D26").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E26").Value = "2105"
$ws.Range("F26").Value = 200000
$ws.Range("G26").Value = 1250000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1143357213"
$ws.Range("D27").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E27").Value = "2104"
$ws.Range("F27").Value = 200000
$ws.Range("G27").Value = 1250000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1143357213"
$ws.Range("D28").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E28").Value = "2103"
$ws.Range("F28").Value = 200000
$ws.Range("G28").Value = 1250000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1143357213"
$ws.Range("D29").Value = "KEVIN ANGULO MARIMON"
$ws.Range("E29").Value = "2102"
$ws.Range("F29").Value = 36341
$ws.Range("G29").Value = 1250000

# BRENDA LUCIA ALVAREZ PEÑA
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1143414760"
$ws.Range("D30").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E30").Value = "2108"
$ws.Range("F30").Value = 35129
$ws.Range("G30").Value = 2500000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1143414760"
$ws.Range("D31").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E31").Value = "2107"
$ws.Range("F31").Value = 100000
$ws.Range("G31").Value = 2500000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1143414760"
$ws.Range("D32").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E32").Value = "2106"
$ws.Range("F32").Value = 100000
$ws.Range("G32").Value = 2500000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1143414760"
$ws.Range("D33").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E33").Value = "2105"
$ws.Range("F33").Value = 100000
$ws.Range("G33").Value = 2500000

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1143414760"
$ws.Range("D34").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E34").Value = "2104"
$ws.Range("F34").Value = 100000
$ws.Range("G34").Value = 2500000

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1143414760"
$ws.Range("D35").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E35").Value = "2103"
$ws.Range("F35").Value = 100000
$ws.Range("G35").Value = 2500000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1143414760"
$ws.Range("D36").Value = "BRENDA LUCIA ALVAREZ PEÑA"
$ws.Range("E36").Value = "2102"
$ws.Range("F36").Value = 36341
$ws.Range("G36").Value = 2500000
